$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Insert a new column at U (21st column). Copying column T first means the
# newly inserted column inherits T's formatting/style (header cell style).
$ws.Columns("T:T").Copy()
$ws.Columns("U:U").Insert()

# New header text for the inserted column
$ws.Range("U1").Value = "Discount Group Consumables"

# Match the new column's width as closely as this host's column-width model allows
$ws.Columns("U").ColumnWidth = $ws.Columns("T").ColumnWidth

# Refresh the AutoFilter so it covers the new last column (A1:AT1). Calling
# AutoFilter() again while filtering is already on toggles it off, so turn
# it off first, then re-enable it with the correct range.
$ws.Range("A1:AS1").AutoFilter() | Out-Null
$ws.Range("A1:AT1").AutoFilter() | Out-Null

# Keep the _FilterDatabase defined name in sync with the new filter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$AT`$1"
    }
}

# Restore the user's selection on the sheet
$ws.Range("U5").Select()

Write-Output "done"
